$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    $cell.Value = $old -replace '^REPSWITCH1_Practice/', 'Pictures_Practice/'
}
